# Populate the newly-added enforcement-action rows (38-48) with their
# Outcome / Cause of Action / Civil-or-Criminal / Token / Project Name /
# Blockchain / Amount / Securities Act / Exchange Act / SEC Office data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 38
$ws.Cells.Item(38, 4).Value = "Ongoing"
$ws.Cells.Item(38, 5).Value = "Unregistered Broker-Dealer"
$ws.Cells.Item(38, 6).Value = "Civil"
$ws.Cells.Item(38, 7).Value = "N/A"
$ws.Cells.Item(38, 8).Value = "AirBit Club"
$ws.Cells.Item(38, 9).Value = "N/A"
$ws.Cells.Item(38, 10).Value = 705000
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 12).Value = 1
$ws.Cells.Item(38, 13).Value = "Southern New York"

# Row 39
$ws.Cells.Item(39, 4).Value = "Ongoing"
$ws.Cells.Item(39, 5).Value = "Unregistered Offering and Fraud"
$ws.Cells.Item(39, 6).Value = "Civil"
$ws.Cells.Item(39, 7).Value = "Boon Coins"
$ws.Cells.Item(39, 8).Value = "Boon.Tech"
$ws.Cells.Item(39, 9).Value = "Ethereum"
$ws.Cells.Item(39, 10).Value = 5000000
$ws.Cells.Item(39, 11).Value = 1
$ws.Cells.Item(39, 12).Value = 1
$ws.Cells.Item(39, 13).Value = "San Francisco"

# Row 40
$ws.Cells.Item(40, 4).Value = "Ongoing"
$ws.Cells.Item(40, 5).Value = "Unregistered Offering and Fraud"
$ws.Cells.Item(40, 6).Value = "Civil and Criminal"
$ws.Cells.Item(40, 7).Value = "ABTC"
$ws.Cells.Item(40, 8).Value = "NAC Foundation"
$ws.Cells.Item(40, 9).Value = "N/A"
$ws.Cells.Item(40, 10).Value = 5600000
$ws.Cells.Item(40, 11).Value = 1
$ws.Cells.Item(40, 12).Value = 1
$ws.Cells.Item(40, 13).Value = "Northern California"

# Row 41
$ws.Cells.Item(41, 4).Value = "Ongoing"
$ws.Cells.Item(41, 5).Value = "Fraud"
$ws.Cells.Item(41, 6).Value = "Civil"
$ws.Cells.Item(41, 7).Value = "N/A"
$ws.Cells.Item(41, 8).Value = "Hvizdzak Capital Management LLC"
$ws.Cells.Item(41, 9).Value = "N/A"
$ws.Cells.Item(41, 10).Value = 31000000
$ws.Cells.Item(41, 11).Value = 1
$ws.Cells.Item(41, 12).Value = 1
$ws.Cells.Item(41, 13).Value = "Chicago"

# Row 42
$ws.Cells.Item(42, 4).Value = "Settlement"
$ws.Cells.Item(42, 5).Value = "Unregistered Offering"
$ws.Cells.Item(42, 6).Value = "Civil "
$ws.Cells.Item(42, 7).Value = "CAT"
$ws.Cells.Item(42, 8).Value = "BitClave"
$ws.Cells.Item(42, 9).Value = "Ethereum"
$ws.Cells.Item(42, 10).Value = 25500000
$ws.Cells.Item(42, 11).Value = 1
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = "San Francisco"

# Row 43
$ws.Cells.Item(43, 4).Value = "Ongoing"
$ws.Cells.Item(43, 5).Value = "Unregistered Offering and Fraud"
$ws.Cells.Item(43, 6).Value = "Civil"
$ws.Cells.Item(43, 7).Value = "BTC"
$ws.Cells.Item(43, 8).Value = "MMT Distributions, LLC"
$ws.Cells.Item(43, 9).Value = "Bitcoin"
$ws.Cells.Item(43, 10).Value = 12000000
$ws.Cells.Item(43, 11).Value = 1
$ws.Cells.Item(43, 12).Value = 1
$ws.Cells.Item(43, 13).Value = "Salt Lake"

# Row 44
$ws.Cells.Item(44, 4).Value = "Ongoing"
$ws.Cells.Item(44, 5).Value = "Unregistered Offering"
$ws.Cells.Item(44, 6).Value = "Civil"
$ws.Cells.Item(44, 7).Value = "DROP"
$ws.Cells.Item(44, 8).Value = "Dropil, Inc."
$ws.Cells.Item(44, 9).Value = "N/A"
$ws.Cells.Item(44, 10).Value = 1800000
$ws.Cells.Item(44, 11).Value = 1
$ws.Cells.Item(44, 12).Value = 1
$ws.Cells.Item(44, 13).Value = "Los Angeles"

# Row 45
$ws.Cells.Item(45, 4).Value = "Ongoing"
$ws.Cells.Item(45, 5).Value = "Unregistered Offering and Fraud"
$ws.Cells.Item(45, 6).Value = "Civil"
$ws.Cells.Item(45, 7).Value = "Meta 1"
$ws.Cells.Item(45, 8).Value = "Pramana Capital Inc."
$ws.Cells.Item(45, 9).Value = "N/A"
$ws.Cells.Item(45, 10).Value = 4300000
$ws.Cells.Item(45, 11).Value = 1
$ws.Cells.Item(45, 12).Value = 1
$ws.Cells.Item(45, 13).Value = "Western District of Texas"

# Row 46
$ws.Cells.Item(46, 4).Value = "Settlement"
$ws.Cells.Item(46, 5).Value = "Anti-touting"
$ws.Cells.Item(46, 6).Value = "Civil"
$ws.Cells.Item(46, 7).Value = "B2G"
$ws.Cells.Item(46, 8).Value = "Steven Seagal"
$ws.Cells.Item(46, 9).Value = "N/A"
$ws.Cells.Item(46, 10).Value = 1000000
$ws.Cells.Item(46, 11).Value = 1
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = "New York"

# Row 47
$ws.Cells.Item(47, 4).Value = "Settlement"
$ws.Cells.Item(47, 5).Value = "Unregistered Offering"
$ws.Cells.Item(47, 6).Value = "Civil"
$ws.Cells.Item(47, 7).Value = "ENG"
$ws.Cells.Item(47, 8).Value = "Enigma MPC"
$ws.Cells.Item(47, 9).Value = "Ethereum"
$ws.Cells.Item(47, 10).Value = 45000000
$ws.Cells.Item(47, 11).Value = 1
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = "Boston"

# Row 48
$ws.Cells.Item(48, 4).Value = "Ongoing"
$ws.Cells.Item(48, 5).Value = "Unregistered Offering and Fraud"
$ws.Cells.Item(48, 6).Value = "Civil"
$ws.Cells.Item(48, 7).Value = "N/A"
$ws.Cells.Item(48, 8).Value = "Q3 I"
$ws.Cells.Item(48, 9).Value = "N/A"
$ws.Cells.Item(48, 10).Value = 33000000
$ws.Cells.Item(48, 11).Value = 1
$ws.Cells.Item(48, 12).Value = 1
$ws.Cells.Item(48, 13).Value = "Miami"

# Restore the author's final on-screen selection/scroll state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C48").Select()
